# Update countries & provincias Spain
# Refreshes the COVID case-count snapshot: new "last updated" timestamp,
# revised totals for a handful of countries, and a reorder of two
# tied (13 total-case) small territories -- Montserrat now outranks
# Islas Malvinas, so their data rows swap places in the sorted list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "last updated" timestamp -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Octubre de 2020 a las 06:57"

# --- India (row 5) ----------------------------------------------------------
$ws.Range("B5").Value = 7761312
$ws.Range("C5").Value = 1672
$ws.Range("D5").Value = 6948497
$ws.Range("E5").Value = 695479

# --- Tailandia (row 149) -----------------------------------------------------
$ws.Range("B149").Value = 3727
$ws.Range("C149").Value = 8
$ws.Range("D149").Value = 3518
$ws.Range("E149").Value = 150

# --- Nueva Zelanda (row 164) --------------------------------------------------
$ws.Range("B164").Value = 1923
$ws.Range("C164").Value = 9
$ws.Range("D164").Value = 1832
$ws.Range("E164").Value = 66

# --- Butan (row 186) ----------------------------------------------------------
$ws.Range("B186").Value = 336
$ws.Range("C186").Value = 4
$ws.Range("E186").Value = 30

# --- Montserrat / Islas Malvinas swap places (rows 216 & 217) ---------------
# Row 216 used to be Islas Malvinas, row 217 used to be Montserrat.
# After the refresh Montserrat (13,0,12,0,0,0,1) sorts ahead of
# Islas Malvinas (13,0,13,0,0,0,0), so the two rows trade both their
# country labels and their data.
$ws.Range("A216").Value = "Montserrat"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 1

$ws.Range("A217").Value = "Islas Malvinas"
$ws.Range("B217").Value = 13
$ws.Range("C217").Value = 0
$ws.Range("D217").Value = 13
$ws.Range("E217").Value = 0
$ws.Range("F217").Value = 0
$ws.Range("G217").Value = 0
$ws.Range("H217").Value = 0

Write-Host "Update complete"
